$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 5) with the same shape as existing rows 2-4:
# Name = "Ayush", Mobile Number = 8368547181, Status = "fraud"
$ws.Range("A5").Value = "Ayush"
$ws.Range("B5").Value = 8368547181
$ws.Range("C5").Value = "fraud"
